# Update cryptos list: price (D) and 1h volume change (E) columns for latest data pull.
# Some coins (Maker/ImmutableX, Aave/Frax, Algorand/USDD) swapped ranking rows.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "25.915.34"
$ws.Range("E2").Value = "  +0.25%  "

$ws.Range("D3").Value = "1.646.31"
$ws.Range("E3").Value = "  +0.63%  "

$ws.Range("E4").Value = "  +0.42%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "215.58"
$ws.Range("E5").Value = "  +0.12%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.5106"
$ws.Range("E6").Value = "  +1.65%  "

$ws.Range("E8").Value = "  +0.09%  "

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.06419"
$ws.Range("E9").Value = "  +0.13%  "

$ws.Range("E10").Value = "  +0.40%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.07760"
$ws.Range("E11").Value = "  +0.70%  "

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "4.312"
$ws.Range("E12").Value = "  +1.60%  "

$ws.Range("D13").Value = "1.654.14"
$ws.Range("E13").Value = "  +1.05%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "0.5480"
$ws.Range("E14").Value = "  +0.59%  "

$ws.Range("D15").Value = "0.0₅7903"
$ws.Range("E15").Value = "  -0.41%  "

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "64.94"
$ws.Range("E16").Value = "  +2.30%  "

$ws.Range("D17").Value = "25.993.04"
$ws.Range("E17").Value = "  +0.50%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "197.74"
$ws.Range("E19").Value = "  -2.50%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "4.437"
$ws.Range("E20").Value = "  +2.58%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "10.04"
$ws.Range("E21").Value = "  +0.88%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "6.061"
$ws.Range("E22").Value = "  +1.37%  "

$ws.Range("E23").Value = "  +0.45%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "1.854"
$ws.Range("E24").Value = "  -3.63%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "140.42"
$ws.Range("E25").Value = "  -0.56%  "

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "0.1148"
$ws.Range("E26").Value = "  +0.32%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "6.897"
$ws.Range("E27").Value = "  +2.81%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "15.77"
$ws.Range("E28").Value = "  +0.53%  "

$ws.Range("E29").Value = "  -0.11%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "0.05011"
$ws.Range("E30").Value = "  +0.14%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "3.276"
$ws.Range("E31").Value = "  +0.44%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "3.203"
$ws.Range("E32").Value = "  +0.86%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "1.544"
$ws.Range("E33").Value = "  +0.61%  "

$ws.Range("E34").Value = "  +0.16%  "

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.8951"
$ws.Range("E35").Value = "  +0.11%  "

$ws.Range("E36").Value = "  -0.69%  "

$ws.Range("B37").Value = "Maker"
$ws.Range("C37").Value = "https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr"
$ws.Range("D37").Value = "1.133.85"
$ws.Range("E37").Value = "  -3.59%  "

$ws.Range("B38").Value = "ImmutableX"
$ws.Range("C38").Value = "https://coinranking.com/coin/Z96jIvLU7+immutablex-imx"
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.5535"
$ws.Range("E38").Value = "  -1.44%  "

$ws.Range("E39").Value = "  +0.45%  "

$ws.Range("E40").Value = "  +0.43%  "

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "5.660"
$ws.Range("E41").Value = "  -0.26%  "

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.8146"
$ws.Range("E42").Value = "  +0.95%  "

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "99.66"
$ws.Range("E43").Value = "  +0.27%  "

$ws.Range("D44").Value = "0.0₈124"
$ws.Range("E44").Value = "  +8.08%  "

$ws.Range("D45").Value = "1.785.11"
$ws.Range("E45").Value = "  +0.69%  "

$ws.Range("E46").Value = "  +0.40%  "

$ws.Range("B47").Value = "Aave"
$ws.Range("C47").Value = "https://coinranking.com/coin/ixgUfzmLR+aave-aave"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "55.34"
$ws.Range("E47").Value = "  +1.06%  "

$ws.Range("B48").Value = "Frax"
$ws.Range("C48").Value = "https://coinranking.com/coin/KfWtaeV1W+frax-frax"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "1.006"
$ws.Range("E48").Value = "  +0.17%  "

$ws.Range("E49").Value = "  +0.60%  "

$ws.Range("B50").Value = "Algorand"
$ws.Range("C50").Value = "https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.09565"
$ws.Range("E50").Value = "  +2.92%  "

$ws.Range("B51").Value = "USDD"
$ws.Range("C51").Value = "https://coinranking.com/coin/z2PZIKQL7+usdd-usdd"
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "1.006"
$ws.Range("E51").Value = "  +0.07%  "
